$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 17966.666
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 17966.666
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 17966.666
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -18424.666

# Row 132
$ws.Range("H132").Value = 2741.68
$ws.Range("I132").Value = 1096.4
$ws.Range("J132").Value = 9322.799999999999
$ws.Range("K132").Value = 3289.2
$ws.Range("L132").Value = 27968.4
$ws.Range("M132").Value = -759.2000000000003
$ws.Range("N132").Value = -33028.39999999999

# Row 137
$ws.Range("H137").Value = 779.4878
$ws.Range("I137").Value = 654.6667
$ws.Range("J137").Value = 877.1739
$ws.Range("K137").Value = 1964.0001
$ws.Range("L137").Value = 2631.5217
$ws.Range("M137").Value = 585.9999
$ws.Range("N137").Value = -7731.5217

# Row 138
$ws.Range("H138").Value = 2701.6196
$ws.Range("I138").Value = 1242.3489
$ws.Range("J138").Value = 4942.643
$ws.Range("K138").Value = 3727.0467
$ws.Range("L138").Value = 14827.929
$ws.Range("M138").Value = 1412.9533
$ws.Range("N138").Value = -25107.929

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 342948.34
$ws.Range("I32").Value = 2345.81
$ws.Range("J32").Value = 5724468.5
$ws.Range("K32").Value = 2345.81
$ws.Range("L32").Value = 5724468.5
$ws.Range("M32").Value = -2058.81
$ws.Range("N32").Value = -5725042.5

# Row 61
$ws.Range("H61").Value = 1007.614
$ws.Range("I61").Value = 654.7805
$ws.Range("J61").Value = 1911.75
$ws.Range("K61").Value = 654.7805
$ws.Range("L61").Value = 1911.75
$ws.Range("M61").Value = -442.7805
$ws.Range("N61").Value = -2335.75

# Row 74
$ws.Range("H74").Value = 147641.78
$ws.Range("I74").Value = 164514.14
$ws.Range("J74").Value = 611.1429000000001
$ws.Range("K74").Value = 164514.14
$ws.Range("L74").Value = 611.1429000000001
$ws.Range("M74").Value = -163640.14
$ws.Range("N74").Value = -2359.1429

# Row 77
$ws.Range("H77").Value = 147641.78
$ws.Range("I77").Value = 164514.14
$ws.Range("J77").Value = 611.1429000000001
$ws.Range("K77").Value = 822570.7000000001
$ws.Range("L77").Value = 3055.7145
$ws.Range("M77").Value = -818202.7000000001
$ws.Range("N77").Value = -11791.7145

# Row 132
$ws.Range("H132").Value = 1279.35
$ws.Range("I132").Value = 960.32434
$ws.Range("J132").Value = 5214
$ws.Range("K132").Value = 2880.97302
$ws.Range("L132").Value = 15642
$ws.Range("M132").Value = -350.9730199999999
$ws.Range("N132").Value = -20702

# Row 136
$ws.Range("H136").Value = 1007.614
$ws.Range("I136").Value = 654.7805
$ws.Range("J136").Value = 1911.75
$ws.Range("K136").Value = 1964.3415
$ws.Range("L136").Value = 5735.25
$ws.Range("M136").Value = 585.6585
$ws.Range("N136").Value = -10835.25

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 30537
$ws.Range("I134").Value = 1950
$ws.Range("J134").Value = 49595
$ws.Range("K134").Value = 5850
$ws.Range("L134").Value = 148785
$ws.Range("M134").Value = -3315
$ws.Range("N134").Value = -153855

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 67.94444
$ws.Range("I7").Value = 87.666664
$ws.Range("J7").Value = 48.22222
$ws.Range("K7").Value = 87.666664
$ws.Range("L7").Value = 48.22222
$ws.Range("M7").Value = 25.333336
$ws.Range("N7").Value = -274.22222

# Row 16
$ws.Range("H16").Value = 1550
$ws.Range("I16").Value = 1552.4445
$ws.Range("J16").Value = 1544.5
$ws.Range("K16").Value = 1552.4445
$ws.Range("L16").Value = 1544.5
$ws.Range("M16").Value = -1265.4445
$ws.Range("N16").Value = -2118.5

# Row 22
$ws.Range("H22").Value = 673.3333
$ws.Range("I22").Value = 611.4286
$ws.Range("J22").Value = 760
$ws.Range("K22").Value = 611.4286
$ws.Range("L22").Value = 760
$ws.Range("M22").Value = -261.4286
$ws.Range("N22").Value = -1460

# Row 68
$ws.Range("H68").Value = 17700.5
$ws.Range("J68").Value = 17700.5
$ws.Range("L68").Value = 17700.5
$ws.Range("N68").Value = -19198.5

# Row 71
$ws.Range("H71").Value = 17700.5
$ws.Range("J71").Value = 17700.5
$ws.Range("L71").Value = 53101.5
$ws.Range("N71").Value = -60589.5

# Row 113
$ws.Range("H113").Value = 1550
$ws.Range("I113").Value = 1552.4445
$ws.Range("J113").Value = 1544.5
$ws.Range("K113").Value = 1552.4445
$ws.Range("L113").Value = 1544.5
$ws.Range("M113").Value = 617.5554999999999
$ws.Range("N113").Value = -5884.5

# Row 132
$ws.Range("H132").Value = 1201.0461
$ws.Range("I132").Value = 1026.5532
$ws.Range("J132").Value = 1656.6666
$ws.Range("K132").Value = 3079.6596
$ws.Range("L132").Value = 4969.9998
$ws.Range("M132").Value = -549.6596
$ws.Range("N132").Value = -10029.9998

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 1514.6666
$ws.Range("I102").Value = 1470.9259
$ws.Range("J102").Value = 1645.8889
$ws.Range("K102").Value = 1470.9259
$ws.Range("L102").Value = 1645.8889
$ws.Range("M102").Value = 151.0741
$ws.Range("N102").Value = -4889.8889

$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 316.66666
$ws.Range("I2").Value = 316.66666
$ws.Range("K2").Value = 316.66666
$ws.Range("M2").Value = -204.66666

# Row 7
$ws.Range("H7").Value = 2804.1035
$ws.Range("I7").Value = 2174.9333
$ws.Range("K7").Value = 2174.9333
$ws.Range("M7").Value = -2062.9333

# Row 46
$ws.Range("H46").Value = 1320.7894
$ws.Range("I46").Value = 1699.3
$ws.Range("J46").Value = 900.2222
$ws.Range("K46").Value = 1699.3
$ws.Range("L46").Value = 900.2222
$ws.Range("M46").Value = -1511.3
$ws.Range("N46").Value = -1276.2222

# Row 126
$ws.Range("H126").Value = 2804.1035
$ws.Range("I126").Value = 2174.9333
$ws.Range("K126").Value = 6524.7999
$ws.Range("M126").Value = -4054.7999

# Row 132
$ws.Range("H132").Value = 1886.1342
$ws.Range("I132").Value = 1950.1666
$ws.Range("J132").Value = 1711.5
$ws.Range("K132").Value = 5850.4998
$ws.Range("L132").Value = 5134.5
$ws.Range("M132").Value = -3320.4998
$ws.Range("N132").Value = -10194.5

# Row 136
$ws.Range("H136").Value = 2454.158
$ws.Range("I136").Value = 1751.9524
$ws.Range("J136").Value = 4420.3335
$ws.Range("K136").Value = 5255.857199999999
$ws.Range("L136").Value = 13261.0005
$ws.Range("M136").Value = -2705.857199999999
$ws.Range("N136").Value = -18361.0005

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 20162330
$ws.Range("I132").Value = 26042646
$ws.Range("J132").Value = 1247.5
$ws.Range("K132").Value = 78127938
$ws.Range("L132").Value = 3742.5
$ws.Range("M132").Value = -78125408
$ws.Range("N132").Value = -8802.5
